$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "gadofosveset trisodium (Lantheus..." - remove the spell-check
#    proofErr markup that wrapped "trisodium" by replacing the
#    surrounding text " trisodium (" (which spans the proofErr
#    boundaries) with itself. Word's Find/Replace merges the whole
#    matched span into a single run and drops the now-superfluous
#    proofErr start/end tags that were fully inside the match.
# -----------------------------------------------------------------
$d.Content.Find.Execute(" trisodium (", $true, $false, $false, $false, $false, $true, 1, $false, " trisodium (", 2) | Out-Null

# -----------------------------------------------------------------
# 2) "Table ?." -> "Table ?" followed by a new run ". " and a
#    (relocated) "_GoBack" bookmark right after it.
# -----------------------------------------------------------------
$word.Selection.HomeKey(6) | Out-Null
$found = $word.Selection.Find.Execute("Table ?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $endPos = $word.Selection.End

    # remove the trailing period from "Table ?."
    $periodRange = $d.Range($endPos, $endPos + 1)
    $periodRange.Delete()

    # the following run currently holds just a single space " ";
    # turn it into ". " (period + space) as its own run
    $spaceRange = $d.Range($endPos, $endPos + 1)
    $spaceRange.Text = ". "

    # Insert a temporary marker character right after the new ". "
    # text. Adding the "_GoBack" bookmark exactly at the end of this
    # paragraph (i.e. immediately before the paragraph mark) is
    # unreliable, so we briefly insert an extra character to move
    # that boundary out of the way, add the bookmark right before
    # it, and then remove the helper character again. Because the
    # bookmark sits before the deleted character it stays put.
    $tempRange = $d.Range($spaceRange.End, $spaceRange.End)
    $tempRange.InsertBefore("X")

    $bmRange = $d.Range($spaceRange.End, $spaceRange.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $bm = $d.Bookmarks("_GoBack")
    $helperRange = $d.Range($bm.End, $bm.End + 1)
    $helperRange.Delete()
}
